$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
Write-Host $ws.Name
